# Fruta / hortaliza, semanal
# Insert a new weekly data row (row 340) into the "Zapallo" sheet, pushing
# all subsequent rows (old 340..368) down by one to (341..369).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 340; existing rows 340-368 shift to 341-369.
$ws.Rows(340).Insert()

# Populate the newly inserted row 340 with this week's record.
$ws.Range("A340").Value = 5
$ws.Range("B340").Value = "Macroferia Regional de Talca"
$ws.Range("C340").Value = "Maule"
$ws.Range("D340").Value = 44931
$ws.Range("E340").Value = 7
$ws.Range("F340").Value = 100112045
$ws.Range("G340").Value = "Zapallo"
$ws.Range("H340").Value = "Camote"
$ws.Range("I340").Value = "1a nueva(o)"
$ws.Range("J340").Value = 900
$ws.Range("K340").Value = 500
$ws.Range("L340").Value = 500
$ws.Range("M340").Value = 500
$ws.Range("N340").Value = "$/kilo (volumen en unidades)"
$ws.Range("O340").Value = "Región del Maule"
$ws.Range("P340").Value = 500
$ws.Range("Q340").Value = 1
$ws.Range("R340").Value = "Hortaliza"
